$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 267 (shifts old rows 267-337 down to 268-338)
$ws.Range("A267:R267").EntireRow.Insert()

# Populate the newly inserted row 267 with the new weekly price record
$ws.Cells.Item(267, 1).Value = 4
$ws.Cells.Item(267, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(267, 3).Value = "Los Lagos"
$ws.Cells.Item(267, 4).Value = 44855
$ws.Cells.Item(267, 5).Value = 10
$ws.Cells.Item(267, 6).Value = 100112017
$ws.Cells.Item(267, 7).Value = "Apio"
$ws.Cells.Item(267, 8).Value = "Americana (o)"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 40
$ws.Cells.Item(267, 11).Value = 13000
$ws.Cells.Item(267, 12).Value = 13000
$ws.Cells.Item(267, 13).Value = 13000
$ws.Cells.Item(267, 14).Value = "`$/docena de matas"
$ws.Cells.Item(267, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(267, 16).Value = 2167
$ws.Cells.Item(267, 17).Value = 6
$ws.Cells.Item(267, 18).Value = "Hortaliza"
